$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts existing rows 13:57 down to 14:58,
# extending the data range from A1:T57 to A1:T58).
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with the new weekly record.
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 'Vega Monumental Concepción'
$ws.Range("C13").Value = 'Bíobío'
$ws.Range("D13").Value = 44624
$ws.Range("E13").Value = 8
$ws.Range("F13").Value = 'Fruta'
$ws.Range("G13").Value = 100103
$ws.Range("H13").Value = 'Frutos de hueso (carozo)'
$ws.Range("I13").Value = 100103002
$ws.Range("J13").Value = 'Ciruela'
$ws.Range("K13").Value = 'Black Amber'
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 250
$ws.Range("N13").Value = 9000
$ws.Range("O13").Value = 9500
$ws.Range("P13").Value = 9300
$ws.Range("Q13").Value = '$/bandeja 18 kilos granel'
$ws.Range("R13").Value = "Región de O'Higgins"
$ws.Range("S13").Value = 517
$ws.Range("T13").Value = 18
